$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B2: set value to "Never" and style it like the other bordered/dated
# cells in the table (date number format, bold header font, centered,
# bright-blue fill) to flag that this car never had an oil change logged.
$rng = $ws.Range("B2")
$rng.Value = "Never"
$rng.NumberFormat = "m/d/yyyy;@"
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.Interior.Color = 15773696     # RGB(0, 176, 240) -> BGR-encoded OLE color

# Move/save the active selection to B2, matching the edited cell.
$rng.Select()
